$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new column header "Unit Tested?" in J2, matching the style of the
# other header cells (row 2 header style).
$ws.Range("J2").Value = "Unit Tested?"
$ws.Range("J2").Style = $ws.Range("I2").Style

# Set the new column's width to match what Excel auto-sized it to.
$ws.Columns.Item(10).ColumnWidth = 12.7109375

# Move the selection / view the same way the author's session ended up:
# scrolled right so column H is the first visible column, with J2 selected.
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("J2").Select()
